$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section ("LOB1004: ...")
# and the paragraph holding the site footer copyright notice. Everything
# between them (an empty paragraph, the "Ver no Jupiter ..." paragraph and
# the "(c) 2020 ..." paragraph itself) is scraped boilerplate that was
# removed from the rebuilt site page, so delete that whole span.
$reqPara = $null
$copyrightPara = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOB1004: C*lculo II (Requisito fraco)*") {
        $reqPara = $p
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $copyrightPara = $p
    }
}

if ($reqPara -ne $null -and $copyrightPara -ne $null) {
    $startPara = $reqPara.Next()
    $start = $startPara.Range.Start
    $end = $copyrightPara.Range.End
    $killRange = $d.Range($start, $end)
    $killRange.Delete()
}
